$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value2 = 1001
$ws.Cells.Item(2,2).Value2 = 18
$ws.Cells.Item(2,3).Value2 = 30
$ws.Cells.Item(2,4).Value2 = 75
$ws.Cells.Item(2,5).Value2 = 60
$ws.Cells.Item(2,6).Value2 = 72

$ws.Cells.Item(3,1).Value2 = 501
$ws.Cells.Item(3,2).Value2 = 9
$ws.Cells.Item(3,3).Value2 = 52
$ws.Cells.Item(3,4).Value2 = 30
$ws.Cells.Item(3,5).Value2 = 75
$ws.Cells.Item(3,6).Value2 = 45

$ws.Cells.Item(4,1).Value2 = 1203
$ws.Cells.Item(4,2).Value2 = 3
$ws.Cells.Item(4,3).Value2 = 15
$ws.Cells.Item(4,4).Value2 = 15
$ws.Cells.Item(4,5).Value2 = 15
$ws.Cells.Item(4,6).Value2 = 15

$ws.Cells.Item(5,1).Value2 = 902
$ws.Cells.Item(5,2).Value2 = 1
$ws.Cells.Item(5,3).Value2 = 0
$ws.Cells.Item(5,4).Value2 = 0
$ws.Cells.Item(5,5).Value2 = 0
$ws.Cells.Item(5,6).Value2 = 0

$ws.Cells.Item(6,1).Value2 = 601
$ws.Cells.Item(6,2).Value2 = 9
$ws.Cells.Item(6,3).Value2 = 60
$ws.Cells.Item(6,4).Value2 = 67
$ws.Cells.Item(6,5).Value2 = 60
$ws.Cells.Item(6,6).Value2 = 42

$ws.Cells.Item(7,1).Value2 = 801
$ws.Cells.Item(7,2).Value2 = 3
$ws.Cells.Item(7,3).Value2 = 67
$ws.Cells.Item(7,4).Value2 = 65
$ws.Cells.Item(7,5).Value2 = 52
$ws.Cells.Item(7,6).Value2 = 45

$ws.Cells.Item(8,1).Value2 = 101
$ws.Cells.Item(8,2).Value2 = 9
$ws.Cells.Item(8,3).Value2 = 30
$ws.Cells.Item(8,4).Value2 = 15
$ws.Cells.Item(8,5).Value2 = 60
$ws.Cells.Item(8,6).Value2 = 15

$ws.Cells.Item(9,1).Value2 = 901
$ws.Cells.Item(9,2).Value2 = 16
$ws.Cells.Item(9,3).Value2 = 15
$ws.Cells.Item(9,4).Value2 = 45
$ws.Cells.Item(9,5).Value2 = 60
$ws.Cells.Item(9,6).Value2 = 60

$ws.Cells.Item(10,1).Value2 = 401
$ws.Cells.Item(10,2).Value2 = 9
$ws.Cells.Item(10,3).Value2 = 48
$ws.Cells.Item(10,4).Value2 = 67
$ws.Cells.Item(10,5).Value2 = 75
$ws.Cells.Item(10,6).Value2 = 45

$ws.Cells.Item(11,1).Value2 = 701
$ws.Cells.Item(11,2).Value2 = 3
$ws.Cells.Item(11,3).Value2 = 90
$ws.Cells.Item(11,4).Value2 = 45
$ws.Cells.Item(11,5).Value2 = 97
$ws.Cells.Item(11,6).Value2 = 15

$ws.Cells.Item(12,1).Value2 = 1201
$ws.Cells.Item(12,2).Value2 = 2
$ws.Cells.Item(12,3).Value2 = 10
$ws.Cells.Item(12,4).Value2 = 10
$ws.Cells.Item(12,5).Value2 = 10
$ws.Cells.Item(12,6).Value2 = 10

$ws.Cells.Item(13,1).Value2 = 1202
$ws.Cells.Item(13,2).Value2 = 2
$ws.Cells.Item(13,3).Value2 = 10
$ws.Cells.Item(13,4).Value2 = 10
$ws.Cells.Item(13,5).Value2 = 10
$ws.Cells.Item(13,6).Value2 = 10

$ws.Cells.Item(14,1).Value2 = 301
$ws.Cells.Item(14,2).Value2 = 6
$ws.Cells.Item(14,3).Value2 = 45
$ws.Cells.Item(14,4).Value2 = 30
$ws.Cells.Item(14,5).Value2 = 60
$ws.Cells.Item(14,6).Value2 = 45

$ws.Cells.Item(15,1).Value2 = 201
$ws.Cells.Item(15,2).Value2 = 9
$ws.Cells.Item(15,3).Value2 = 30
$ws.Cells.Item(15,4).Value2 = 15
$ws.Cells.Item(15,5).Value2 = 45
$ws.Cells.Item(15,6).Value2 = 30

$ws.Cells.Item(16,1).Value2 = 2
$ws.Cells.Item(16,2).Value2 = 0
$ws.Cells.Item(16,3).Value2 = 2
$ws.Cells.Item(16,4).Value2 = 2
$ws.Cells.Item(16,5).Value2 = 2
$ws.Cells.Item(16,6).Value2 = 2

$ws.Cells.Item(17,1).Value2 = 802
$ws.Cells.Item(17,2).Value2 = 0
$ws.Cells.Item(17,3).Value2 = 4
$ws.Cells.Item(17,4).Value2 = 5
$ws.Cells.Item(17,5).Value2 = 4
$ws.Cells.Item(17,6).Value2 = 0

$ws.Cells.Item(18,1).Value2 = 3
$ws.Cells.Item(18,2).Value2 = 0
$ws.Cells.Item(18,3).Value2 = 3
$ws.Cells.Item(18,4).Value2 = 3
$ws.Cells.Item(18,5).Value2 = 3
$ws.Cells.Item(18,6).Value2 = 3

$ws.Cells.Item(19,1).Value2 = 1101
$ws.Cells.Item(19,2).Value2 = 0
$ws.Cells.Item(19,3).Value2 = 15
$ws.Cells.Item(19,4).Value2 = 30
$ws.Cells.Item(19,5).Value2 = 30
$ws.Cells.Item(19,6).Value2 = 0

$ws.Cells.Item(20,1).Value2 = 1
$ws.Cells.Item(20,2).Value2 = 0
$ws.Cells.Item(20,3).Value2 = 2
$ws.Cells.Item(20,4).Value2 = 2
$ws.Cells.Item(20,5).Value2 = 2
$ws.Cells.Item(20,6).Value2 = 2

$ws.Cells.Item(21,1).Value2 = 502
$ws.Cells.Item(21,2).Value2 = 0
$ws.Cells.Item(21,3).Value2 = 4
$ws.Cells.Item(21,4).Value2 = 0
$ws.Cells.Item(21,5).Value2 = 0
$ws.Cells.Item(21,6).Value2 = 0

